$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, pushing the existing rows 68-81 down to 69-82.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly record.
$ws.Cells.Item(68, 1).Value = 9
$ws.Cells.Item(68, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(68, 3).Value = "Metropolitana"
$ws.Cells.Item(68, 4).Value = 44841
$ws.Cells.Item(68, 5).Value = 13
$ws.Cells.Item(68, 6).Value = 100112029
$ws.Cells.Item(68, 7).Value = "Orégano"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 16
$ws.Cells.Item(68, 11).Value = 18000
$ws.Cells.Item(68, 12).Value = 18000
$ws.Cells.Item(68, 13).Value = 18000
$ws.Cells.Item(68, 14).Value = "$/docena de atados"
$ws.Cells.Item(68, 15).Value = "Región Metropolitana"
$ws.Cells.Item(68, 16).Value = 6000
$ws.Cells.Item(68, 17).Value = 3
$ws.Cells.Item(68, 18).Value = "Hortaliza"
